# Budget list update: fill in rows 17-22 (Item/Service Name, Quantity,
# Shipping Fee, Price per Unit) on the "Alinanlar" sheet with the new
# components that were ordered, and move the on-screen selection down to
# where the new data was entered (A23, the first still-empty row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: L298N
$ws.Range("A17").Value = "L298N"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0.44

# Row 18: TB6612FNG (trailing space kept, as in the source sheet)
$ws.Range("A18").Value = "TB6612FNG "
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1.09

# Row 19: 0.91" white/Blue OLED
$ws.Range("A19").Value = "0.91`" white/Blue OLED "
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 1.56

# Row 20: GY-273 HMC5883L (leading space + double inner space kept)
$ws.Range("A20").Value = " GY-273  HMC5883L"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 1.82

# Row 21: MEGA2560
$ws.Range("A21").Value = "MEGA2560"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 2.11
$ws.Range("D21").Value = 7.95

# Row 22: TP5100
$ws.Range("A22").Value = "TP5100"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0.61

# Reflect the author's on-screen selection after entering the data.
$ws.Range("A23").Select()
